$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target roster data (player, position, team) for rows 2-19
$data = @(
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Spencer Dinwiddie", "PG,SG", "Dallas Mavericks"),
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Trayce Jackson-Davis", "PF,C", "Golden State Warriors"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
